$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1282.1052
$ws.Range("I19").Value = 1292.3125
$ws.Range("K19").Value = 1292.3125
$ws.Range("M19").Value = -1117.3125
$ws.Range("H40").Value = 3841.8462
$ws.Range("I40").Value = 2999.3333
$ws.Range("J40").Value = 3951.739
$ws.Range("K40").Value = 2999.3333
$ws.Range("L40").Value = 3951.739
$ws.Range("M40").Value = -2824.3333
$ws.Range("N40").Value = -4301.739
$ws.Range("H86").Value = 361115360
$ws.Range("I86").Value = 400003870
$ws.Range("K86").Value = 400003870
$ws.Range("M86").Value = -400002747
$ws.Range("H89").Value = 361115360
$ws.Range("I89").Value = 400003870
$ws.Range("K89").Value = 2000019350
$ws.Range("M89").Value = -2000013734
$ws.Range("H116").Value = 35723024
$ws.Range("I116").Value = 62509176
$ws.Range("J116").Value = 8159.3335
$ws.Range("K116").Value = 62509176
$ws.Range("L116").Value = 8159.3335
$ws.Range("M116").Value = -62505734
$ws.Range("N116").Value = -15043.3335
$ws.Range("H132").Value = 3480.84
$ws.Range("I132").Value = 3210.6191
$ws.Range("J132").Value = 4899.5
$ws.Range("K132").Value = 9631.8573
$ws.Range("L132").Value = 14698.5
$ws.Range("M132").Value = -7101.8573
$ws.Range("N132").Value = -19758.5
$ws.Range("H135").Value = 2035.7059
$ws.Range("I135").Value = 1535.909
$ws.Range("J135").Value = 2952
$ws.Range("K135").Value = 13823.181
$ws.Range("L135").Value = 26568
$ws.Range("M135").Value = -11288.181
$ws.Range("N135").Value = -31638
$ws.Range("H138").Value = 2339.1226
$ws.Range("I138").Value = 943.35
$ws.Range("J138").Value = 2697.013
$ws.Range("K138").Value = 2830.05
$ws.Range("L138").Value = 8091.039
$ws.Range("M138").Value = 2309.95
$ws.Range("N138").Value = -18371.039
$ws.Range("H141").Value = 9246.241
$ws.Range("I141").Value = 8623.799999999999
$ws.Range("J141").Value = 9913.143
$ws.Range("K141").Value = 25871.4
$ws.Range("L141").Value = 29739.429
$ws.Range("M141").Value = -20691.4
$ws.Range("N141").Value = -40099.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17490148
$ws.Range("I32").Value = 20725440
$ws.Range("K32").Value = 20725440
$ws.Range("M32").Value = -20725153
$ws.Range("H74").Value = 2385.976
$ws.Range("I74").Value = 2097.5278
$ws.Range("K74").Value = 2097.5278
$ws.Range("M74").Value = -1223.5278
$ws.Range("H77").Value = 2385.976
$ws.Range("I77").Value = 2097.5278
$ws.Range("K77").Value = 10487.639
$ws.Range("M77").Value = -6119.638999999999
$ws.Range("H110").Value = 1448.8235
$ws.Range("I110").Value = 1320.625
$ws.Range("J110").Value = 3500
$ws.Range("K110").Value = 1320.625
$ws.Range("L110").Value = 3500
$ws.Range("M110").Value = 724.375
$ws.Range("N110").Value = -7590
$ws.Range("H132").Value = 5029.5557
$ws.Range("I132").Value = 5187.4287
$ws.Range("J132").Value = 4477
$ws.Range("K132").Value = 15562.2861
$ws.Range("L132").Value = 13431
$ws.Range("M132").Value = -13032.2861
$ws.Range("N132").Value = -18491

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2513374.2
$ws.Range("I7").Value = 5004999.5
$ws.Range("J7").Value = 21749
$ws.Range("K7").Value = 5004999.5
$ws.Range("L7").Value = 21749
$ws.Range("M7").Value = -5004886.5
$ws.Range("N7").Value = -21975
$ws.Range("H99").Value = 2823.739
$ws.Range("I99").Value = 2149.7693
$ws.Range("K99").Value = 2149.7693
$ws.Range("M99").Value = -651.7692999999999
$ws.Range("H134").Value = 2859424
$ws.Range("I134").Value = 3403147.8
$ws.Range("J134").Value = 4873.75
$ws.Range("K134").Value = 10209443.4
$ws.Range("L134").Value = 14621.25
$ws.Range("M134").Value = -10206908.4
$ws.Range("N134").Value = -19691.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 50000
$ws.Range("J50").Value = 50000
$ws.Range("L50").Value = 50000
$ws.Range("H86").Value = 24336.408
$ws.Range("I86").Value = 29098.6
$ws.Range("K86").Value = 29098.6
$ws.Range("M86").Value = -27975.6
$ws.Range("H89").Value = 24336.408
$ws.Range("I89").Value = 29098.6
$ws.Range("K89").Value = 145493
$ws.Range("M89").Value = -139877
$ws.Range("H132").Value = 3407.6572
$ws.Range("I132").Value = 3032.1072
$ws.Range("K132").Value = 9096.321599999999
$ws.Range("M132").Value = -6566.321599999999
$ws.Range("H134").Value = 1888.8
$ws.Range("I134").Value = 1891.862
$ws.Range("K134").Value = 5675.586
$ws.Range("M134").Value = -3140.586
$ws.Range("N50").Value = -51250

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 229
$ws.Range("I2").Value = 142.66667
$ws.Range("J2").Value = 488
$ws.Range("K2").Value = 856.0000200000001
$ws.Range("L2").Value = 2928
$ws.Range("M2").Value = -743.0000200000001
$ws.Range("N2").Value = -3154
$ws.Range("H12").Value = 110.111115
$ws.Range("J12").Value = 110.111115
$ws.Range("L12").Value = 330.333345
$ws.Range("N12").Value = -676.333345
$ws.Range("H14").Value = 776.2222
$ws.Range("I14").Value = 776.2222
$ws.Range("K14").Value = 2328.6666
$ws.Range("M14").Value = -2155.6666
$ws.Range("H38").Value = 58.733334
$ws.Range("I38").Value = 51.3
$ws.Range("J38").Value = 73.59999999999999
$ws.Range("K38").Value = 153.9
$ws.Range("L38").Value = 220.8
$ws.Range("M38").Value = 193.1
$ws.Range("N38").Value = -914.8
$ws.Range("H92").Value = 5747.5
$ws.Range("I92").Value = 4997.5
$ws.Range("J92").Value = 6497.5
$ws.Range("K92").Value = 14992.5
$ws.Range("L92").Value = 19492.5
$ws.Range("M92").Value = -13744.5
$ws.Range("N92").Value = -21988.5
$ws.Range("H97").Value = 621.5714
$ws.Range("I97").Value = 554.4
$ws.Range("J97").Value = 789.5
$ws.Range("K97").Value = 1663.2
$ws.Range("L97").Value = 2368.5
$ws.Range("M97").Value = -1167.2
$ws.Range("N97").Value = -3360.5
$ws.Range("H131").Value = 2406.7646
$ws.Range("I131").Value = 1369.8334
$ws.Range("J131").Value = 2972.3635
$ws.Range("K131").Value = 4109.5002
$ws.Range("L131").Value = 8917.0905
$ws.Range("M131").Value = 930.4997999999996
$ws.Range("N131").Value = -18997.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 64349
$ws.Range("J51").Value = 64349
$ws.Range("L51").Value = 64349
$ws.Range("N51").Value = -65367
$ws.Range("H97").Value = 1496.4286
$ws.Range("I97").Value = 719.7
$ws.Range("K97").Value = 719.7
$ws.Range("M97").Value = -223.7
$ws.Range("H113").Value = 21604.6
$ws.Range("J113").Value = 26403
$ws.Range("L113").Value = 26403
$ws.Range("N113").Value = -30743
$ws.Range("H122").Value = 1741.7858
$ws.Range("J122").Value = 2151.3333
$ws.Range("L122").Value = 6453.999899999999
$ws.Range("N122").Value = -11353.9999
$ws.Range("H132").Value = 2813.8667
$ws.Range("I132").Value = 2813.8667
$ws.Range("K132").Value = 8441.6001
$ws.Range("M132").Value = -5911.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 30009.5
$ws.Range("I5").Value = 20009
$ws.Range("J5").Value = 40010
$ws.Range("K5").Value = 20009
$ws.Range("L5").Value = 40010
$ws.Range("N5").Value = -40236
$ws.Range("H61").Value = 3003.25
$ws.Range("I61").Value = 2794.3333
$ws.Range("K61").Value = 2794.3333
$ws.Range("M61").Value = -2592.3333
$ws.Range("H82").Value = 2581.5264
$ws.Range("I82").Value = 2474.2727
$ws.Range("K82").Value = 2474.2727
$ws.Range("M82").Value = -2113.2727
$ws.Range("H85").Value = 2581.5264
$ws.Range("I85").Value = 2474.2727
$ws.Range("K85").Value = 2474.2727
$ws.Range("M85").Value = -1226.2727
$ws.Range("H100").Value = 2731
$ws.Range("I100").Value = 2494
$ws.Range("J100").Value = 2849.5
$ws.Range("K100").Value = 2494
$ws.Range("L100").Value = 2849.5
$ws.Range("M100").Value = -1953
$ws.Range("N100").Value = -3931.5
$ws.Range("H113").Value = 3003.25
$ws.Range("I113").Value = 2794.3333
$ws.Range("K113").Value = 2794.3333
$ws.Range("M113").Value = -624.3332999999998
$ws.Range("H136").Value = 12599.7
$ws.Range("I136").Value = 5454.636
$ws.Range("K136").Value = 16363.908
$ws.Range("M136").Value = -13813.908
$ws.Range("M5").Value = -19896

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 62773.5
$ws.Range("J16").Value = 62773.5
$ws.Range("L16").Value = 62773.5
$ws.Range("N16").Value = -63357.5
$ws.Range("H41").Value = 12000
$ws.Range("J41").Value = 12000
$ws.Range("L41").Value = 12000
$ws.Range("H45").Value = 28367.25
$ws.Range("J45").Value = 36323.332
$ws.Range("L45").Value = 36323.332
$ws.Range("N45").Value = -37305.332
$ws.Range("H113").Value = 757
$ws.Range("I113").Value = 592.25
$ws.Range("J113").Value = 866.8333
$ws.Range("K113").Value = 1776.75
$ws.Range("L113").Value = 2600.4999
$ws.Range("M113").Value = 393.25
$ws.Range("N113").Value = -6940.4999
$ws.Range("H132").Value = 2300.9119
$ws.Range("I132").Value = 2294.3572
$ws.Range("K132").Value = 6883.071599999999
$ws.Range("M132").Value = -4353.071599999999
$ws.Range("H136").Value = 8377652.5
$ws.Range("I136").Value = 1670.5
$ws.Range("J136").Value = 14360496
$ws.Range("K136").Value = 5011.5
$ws.Range("L136").Value = 43081488
$ws.Range("M136").Value = -2461.5
$ws.Range("N136").Value = -43086588
$ws.Range("N41").Value = -12780
